# Insert a new data row before current row 45 (shifting existing rows
# 45..129 down to 46..130) and populate the new row with the "Ají" /
# "Americana (o)" observation dated 2023-01-26.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("45:45").Insert()

$ws.Range("A45").Value = 7
$ws.Range("B45").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C45").Value = "Ñuble"
$ws.Range("D45").Value = 44952
$ws.Range("E45").Value = 16
$ws.Range("F45").Value = 100112021
$ws.Range("G45").Value = "Ají"
$ws.Range("H45").Value = "Americana (o)"
$ws.Range("I45").Value = "Primera"
$ws.Range("J45").Value = 30
$ws.Range("K45").Value = 12000
$ws.Range("L45").Value = 12000
$ws.Range("M45").Value = 12000
$ws.Range("N45").Value = "$/caja 15 kilos"
$ws.Range("O45").Value = "Región del Maule"
$ws.Range("P45").Value = 800
$ws.Range("Q45").Value = 15
$ws.Range("R45").Value = "Hortaliza"
